# Applies the "additional scraping" commit:
#  1. Insert a new "Player Info" sheet at the front with the player's bio.
#  2. Rename MATCH_CARD_LINK -> MATCH_CODE on both "ODI Batting" and
#     "ODI Bowling", converting the stored howstat URL into the bare
#     numeric match code.
#  3. Drop the handful of stray empty INNING_NUMBER placeholder cells on
#     "ODI Batting".
#  4. Append a new "ODI Batting Extra" sheet with additional per-match
#     batting detail.
#
# NOTE: worksheet object handles in this host track *positional* slots,
# not stable identity - after any operation that inserts/moves/removes a
# sheet, previously-captured worksheet variables can silently start
# pointing at the wrong tab. So every worksheet reference below is
# (re)fetched by name immediately before it is used.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Player Info" sheet, inserted before "ODI Batting" (becomes sheet 1)
# ---------------------------------------------------------------------
$battingForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingForInsert)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value2 = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$playerInfo.Columns.Item(1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value2 = "3969"
$playerInfo.Cells.Item(2, 2).Value2 = "Glenn James Maxwell"
$playerInfo.Cells.Item(2, 3).Value2 = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value2 = "Right Arm Off Break"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code,
#    and drop the stray empty B-column placeholder cells.
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$battingLastRow = $batting.UsedRange.Rows.Count

$batting.Cells.Item(1, 4).Value2 = "MATCH_CODE"
$batting.Columns.Item(4).NumberFormat = "@"

for ($r = 2; $r -le $battingLastRow; $r++) {
    $linkCell = $batting.Cells.Item($r, 4)
    $link = $linkCell.Value2
    if ($link -ne $null -and $link -ne "") {
        $code = $link.Substring($link.LastIndexOf("=") + 1)
        $linkCell.Value2 = $code
    }

    $innings = $batting.Cells.Item($r, 2)
    if ($innings.Value2 -eq $null -or $innings.Value2 -eq "") {
        $innings.ClearContents()
    }
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare code.
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowlingLastRow = $bowling.UsedRange.Rows.Count

$bowling.Cells.Item(1, 2).Value2 = "MATCH_CODE"
$bowling.Columns.Item(2).NumberFormat = "@"

for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $linkCell = $bowling.Cells.Item($r, 2)
    $link = $linkCell.Value2
    if ($link -ne $null -and $link -ne "") {
        $code = $link.Substring($link.LastIndexOf("=") + 1)
        $linkCell.Value2 = $code
    }
}

# ---------------------------------------------------------------------
# 4. "ODI Batting Extra" sheet, appended after "ODI Bowling" (sheet 4)
# ---------------------------------------------------------------------
$bowlingForInsert = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingForInsert)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value2 = $exHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$extra.Columns.Item(1).NumberFormat = "@"
$extra.Columns.Item(3).NumberFormat = "@"
$extra.Columns.Item(4).NumberFormat = "@"
$extra.Columns.Item(5).NumberFormat = "@"
$extra.Columns.Item(6).NumberFormat = "@"

$exRows = @(
    @("4351", "6", "1", "0", "3.81%", "NO"),
    @("4354", "7", "2", "1", "9.87%", "NO"),
    @("4429", "7", "4", "4", "26.19%", "NO"),
    @("4430", "7", "0", "0", "0.48%", "NO"),
    @("4431", "7", "4", "7", "35.41%", "YES"),
    @("4435", "", "", "", "", "NO"),
    @("4436", "", "", "", "", "NO"),
    @("4437", "7", "3", "4", "20.42%", "NO"),
    @("4594", "", "", "", "", "NO"),
    @("4597", "7", "5", "0", "15.87%", "NO"),
    @("4600", "", "", "", "", "NO"),
    @("4601", "7", "0", "0", "0.39%", "NO"),
    @("4603", "7", "2", "0", "9.76%", "NO"),
    @("4644", "7", "3", "3", "15.92%", "NO"),
    @("4645", "7", "", "", "", "NO"),
    @("4646", "7", "3", "0", "13.48%", "NO"),
    @("4647", "", "", "", "", "NO"),
    @("4648", "7", "1", "1", "12.82%", "NO"),
    @("4649", "6", "0", "1", "5.24%", "NO"),
    @("4725", "7", "1", "0", "4.26%", "NO")
)

$r = 2
foreach ($row in $exRows) {
    $extra.Cells.Item($r, 1).Value2 = $row[0]

    if ($row[1] -ne "") {
        $extra.Cells.Item($r, 2).Value2 = [int]$row[1]
    }
    $extra.Cells.Item($r, 3).Value2 = $row[2]
    $extra.Cells.Item($r, 4).Value2 = $row[3]
    $extra.Cells.Item($r, 5).Value2 = $row[4]
    $extra.Cells.Item($r, 6).Value2 = $row[5]

    $r = $r + 1
}

Write-Output "edit applied"
